# Update odds data for the week-of-2024-11-11 FlashScore sheet.
# The source feed refreshed: most odds columns shifted slightly, a couple of
# pairs of columns got swapped, and the last fixture (old row 7,
# "America De Cali vs Santa Fe") was merged into row 6 (replacing the
# "Once Caldas vs Junior" fixture), so the trailing row is dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Argentinos Jrs vs Banfield
$ws.Range("G2").Value  = 1.75
$ws.Range("H2").Value  = 3.3
$ws.Range("I2").Value  = 5.5
$ws.Range("J2").Value  = 2.5
$ws.Range("K2").Value  = 1.95
$ws.Range("M2").Value  = 1.11
$ws.Range("N2").Value  = 6.5
$ws.Range("X2").Value  = 7
$ws.Range("AG2").Value = 10
$ws.Range("AJ2").Value = 51
$ws.Range("AN2").Value = 3.5
$ws.Range("AO2").Value = 10
$ws.Range("AW2").Value = 6.5
$ws.Range("AZ2").Value = 126
$ws.Range("BD2").Value = 126

# Row 3 - Belgrano vs Instituto
$ws.Range("H3").Value  = 3
$ws.Range("Q3").Value  = 2.5
$ws.Range("R3").Value  = 1.5
$ws.Range("AK3").Value = 29

# Row 4 - Coritiba vs Santos
$ws.Range("G4").Value  = 3.75
$ws.Range("I4").Value  = 2.2
$ws.Range("L4").Value  = 3
$ws.Range("M4").Value  = 1.1
$ws.Range("N4").Value  = 7
$ws.Range("Q4").Value  = 2.4
$ws.Range("R4").Value  = 1.53
$ws.Range("S4").Value  = 1.53
$ws.Range("T4").Value  = 2.38
$ws.Range("AN4").Value = 5.5
$ws.Range("AQ4").Value = 81
$ws.Range("AT4").Value = 2.38

# Row 5 - Paysandu PA vs Brusque
$ws.Range("G5").Value  = 1.53
$ws.Range("H5").Value  = 4
$ws.Range("I5").Value  = 6.5
$ws.Range("J5").Value  = 2.1
$ws.Range("K5").Value  = 2.2
$ws.Range("L5").Value  = 7
$ws.Range("Z5").Value  = 10
$ws.Range("AC5").Value = 8.5
$ws.Range("AD5").Value = 8
$ws.Range("AG5").Value = 13
$ws.Range("AI5").Value = 21
$ws.Range("AN5").Value = 3.25
$ws.Range("AX5").Value = 41
$ws.Range("AZ5").Value = 151
$ws.Range("BA5").Value = 201

# Row 6 - now "America De Cali vs Santa Fe" (was "Once Caldas vs Junior")
$ws.Range("A6").Value  = "ARJPKb8t"
$ws.Range("C6").Value  = "22:30"
$ws.Range("E6").Value  = "America De Cali"
$ws.Range("F6").Value  = "Santa Fe"
$ws.Range("G6").Value  = 1.85
$ws.Range("H6").Value  = 3.1
$ws.Range("I6").Value  = 4.75
$ws.Range("K6").Value  = 1.95
$ws.Range("M6").Value  = 1.1
$ws.Range("N6").Value  = 7
$ws.Range("Q6").Value  = 2.4
$ws.Range("R6").Value  = 1.53
$ws.Range("U6").Value  = 2.2
$ws.Range("V6").Value  = 1.62
$ws.Range("W6").Value  = 5.5
$ws.Range("X6").Value  = 7.5
$ws.Range("Y6").Value  = 9.5
$ws.Range("AA6").Value = 19
$ws.Range("AB6").Value = 41
$ws.Range("AC6").Value = 6.5
$ws.Range("AD6").Value = 6
$ws.Range("AI6").Value = 17
$ws.Range("AM6").Value = 201
$ws.Range("AN6").Value = 3.6
$ws.Range("AS6").Value = 251
$ws.Range("AU6").Value = 9.5
$ws.Range("AV6").Value = 81
$ws.Range("AX6").Value = 29
$ws.Range("BA6").Value = 151

# The old row 7 fixture was folded into row 6 above, so drop the now-duplicate
# trailing row entirely (shrinks the sheet from A1:BD7 to A1:BD6).
$ws.Rows("7:7").Delete()
